$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of data among rows 7, 8 and 10:
#   new row 7  <- old row 10
#   new row 8  <- old row 7
#   new row 10 <- old row 8
# Only columns A (Id), I (Antal), Q (Ost), R (Nord), Z (Starttid) and
# AB (Sluttid) actually carry differing data between these three rows;
# all other columns already hold identical content in rows 7/8/10, so
# rotating just these columns reproduces the full-row rotation.
#
# We stage the row being overwritten first ("row7") in an unused scratch
# row (far below the used range) so its values aren't lost, then shuffle
# the remaining rows, and finally drop the data into row 8 from the
# scratch area before clearing the scratch row again. Using Range.Copy
# (cell to cell) rather than re-typing .Value preserves each cell's
# original data type (numeric vs. text) exactly - e.g. column I holds
# text such as "10", not a number.

$cols = @("A", "I", "Q", "R", "Z", "AB")
$scratchRow = 100

foreach ($col in $cols) {
    $ws.Range("$col$scratchRow").ClearContents()
    $ws.Range("${col}7").Copy($ws.Range("$col$scratchRow"))
}

foreach ($col in $cols) {
    $ws.Range("${col}7").ClearContents()
    $ws.Range("${col}10").Copy($ws.Range("${col}7"))
}

foreach ($col in $cols) {
    $ws.Range("${col}10").ClearContents()
    $ws.Range("${col}8").Copy($ws.Range("${col}10"))
}

foreach ($col in $cols) {
    $ws.Range("${col}8").ClearContents()
    $ws.Range("$col$scratchRow").Copy($ws.Range("${col}8"))
}

foreach ($col in $cols) {
    $ws.Range("$col$scratchRow").ClearContents()
}
